$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.969.87'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.11%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.988.28'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '261.98'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +5.88%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.607'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -2.37%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '55.54'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -5.16%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.372'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -4.55%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0759'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -5.55%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -3.56%  '
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.285.93'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.86%  '
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.06'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -6.20%  '
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.86'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.62%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.764'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -8.89%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.15'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -4.80%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.987.26'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.14%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '36.822.92'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.25%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '69.41'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.14%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0822'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -4.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '232.94'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.93%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.04'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -3.69%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.57'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.17%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.05%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.97'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.85%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.80'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -4.67%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.25'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -3.19%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.127'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -8.04%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.29'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -4.59%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.46%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.50'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -5.99%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0616'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -7.93%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.32'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -5.30%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.39'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.50%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.42'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -3.99%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.55%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.33'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.84%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.92%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.07%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.435.15'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.67%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0909'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -6.56%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0206'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -5.41%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '88.57'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.87%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '15.37'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -6.96%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -3.19%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.90'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.76'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -9.44%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.179.38'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.91%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.91'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -9.03%  '
